$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Row 2 / Row 3 - updated utilization figures
# ---------------------------------------------------------------------
$ws.Range("A2").Value = 0.0017788336845114827
$ws.Range("A3").Value = 0.0017788336845114827
$ws.Range("G3").Value = 462.0
$ws.Range("H3").Value = 142.0
$ws.Range("I3").Value = 3.253520965576172

# ---------------------------------------------------------------------
# 2) Prepare new rows 5, 6 and 7 (row 4 already exists) - copy the
#    cell formatting from the existing sibling rows so number formats /
#    fills / borders / indents stay consistent with the report style.
# ---------------------------------------------------------------------
$ws.Range("B3:J3").Copy() | Out-Null
$ws.Range("B4:J4").PasteSpecial(-4122) | Out-Null
$ws.Range("B5:J5").PasteSpecial(-4122) | Out-Null
$ws.Range("B6:J6").PasteSpecial(-4122) | Out-Null
$ws.Range("B7:J7").PasteSpecial(-4122) | Out-Null

$ws.Range("A4").Copy() | Out-Null
$ws.Range("A5:A7").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Rows 4-7 all sit two levels deep in the outline, same as before for row 4.
$ws.Rows.Item(4).OutlineLevel = 2
$ws.Rows.Item(5).OutlineLevel = 2
$ws.Rows.Item(6).OutlineLevel = 2
$ws.Rows.Item(7).OutlineLevel = 2

# ---------------------------------------------------------------------
# 3) Row 4 - ap_clk_IBUF_BUFG (BUFG)
# ---------------------------------------------------------------------
$ws.Range("A4").Value = 0.0015990337124094367
$ws.Range("B4").Value = "ap_clk_IBUF_BUFG"
$ws.Range("C4").Value = 100.0
$ws.Range("D4").Value = "BUFG"
$ws.Range("E4").Value = "N/A"
$ws.Range("F4").Value = "N/A"
$ws.Range("G4").Value = 460.0
$ws.Range("H4").Value = 140.0
$ws.Range("I4").Value = 3.2857139110565186
$ws.Range("J4").Value = "N/A"

# ---------------------------------------------------------------------
# 4) Row 5 - ap_clk_IBUF_BUFG_inst (BUFG)
# ---------------------------------------------------------------------
$ws.Range("A5").Value = 0.00009299999510403723
$ws.Range("B5").Value = "ap_clk_IBUF_BUFG_inst (BUFG)"
$ws.Range("C5").Value = 100.0
$ws.Range("D5").Value = "BUFG"
$ws.Range("E5").Value = "N/A"
$ws.Range("F5").Value = "N/A"
$ws.Range("G5").Value = 460.0
$ws.Range("H5").Value = 140.0
$ws.Range("I5").Value = 3.2857139110565186
$ws.Range("J5").Value = "Global"

# ---------------------------------------------------------------------
# 5) Row 6 - ap_clk_IBUF
# ---------------------------------------------------------------------
$ws.Range("A6").Value = 0.00008679999882588163
$ws.Range("B6").Value = "ap_clk_IBUF"
$ws.Range("C6").Value = 100.0
$ws.Range("D6").Value = "N/A"
$ws.Range("E6").Value = "N/A"
$ws.Range("F6").Value = "N/A"
$ws.Range("G6").Value = 1.0
$ws.Range("H6").Value = 1.0
$ws.Range("I6").Value = 1.0
$ws.Range("J6").Value = "N/A"

# ---------------------------------------------------------------------
# 6) Row 7 - ap_clk
# ---------------------------------------------------------------------
$ws.Range("A7").Value = 0.0
$ws.Range("B7").Value = "ap_clk"
$ws.Range("C7").Value = 100.0
$ws.Range("D7").Value = "N/A"
$ws.Range("E7").Value = "N/A"
$ws.Range("F7").Value = "N/A"
$ws.Range("G7").Value = 1.0
$ws.Range("H7").Value = 1.0
$ws.Range("I7").Value = 1.0
$ws.Range("J7").Value = "N/A"
